# Refresh cryptocurrency Price (column D) and Volume(1h) (column E)
# values scraped for this run. Numeric-looking Price strings are
# written with a leading apostrophe so Excel keeps them as text
# (matching the original inline-string cell type) instead of
# auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.649.80'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '2.461.02'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('D5').Value = "'573.29"
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').Value = "'146.95"
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -1.42%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').Value = "'29.08"
$ws.Range('E13').Value = '  +1.51%  '
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').Value = '2.903.67'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').Value = '62.656.37'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('D17').Value = '2.464.30'
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = "'10.92"
$ws.Range('E19').Value = '  -1.33%  '
$ws.Range('D20').Value = "'325.13"
$ws.Range('E20').Value = '  -1.52%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').Value = "'2.18"
$ws.Range('E22').Value = '  +2.19%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = "'10.07"
$ws.Range('E24').Value = '  +17.13%  '
$ws.Range('D25').Value = "'65.28"
$ws.Range('E25').Value = '  -1.72%  '
$ws.Range('D26').Value = "'640.52"
$ws.Range('E26').Value = '  -2.79%  '
$ws.Range('D27').Value = '2.581.81'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').Value = '0.0₃0971'
$ws.Range('E28').Value = '  -2.61%  '
$ws.Range('D29').Value = "'0.998"
$ws.Range('E29').Value = '  -17.94%  '
$ws.Range('D30').Value = "'1.43"
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').Value = "'7.91"
$ws.Range('E31').Value = '  -3.59%  '
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('E33').Value = '  -4.23%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('E35').Value = '  +2.46%  '
$ws.Range('E36').Value = '  -0.98%  '
$ws.Range('D37').Value = "'151.71"
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('D38').Value = "'0.368"
$ws.Range('E38').Value = '  -1.56%  '
$ws.Range('D39').Value = "'18.58"
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('D40').Value = "'5.33"
$ws.Range('E40').Value = '  -3.26%  '
$ws.Range('E41').Value = '  +0.24%  '
$ws.Range('D42').Value = "'1.73"
$ws.Range('E42').Value = '  -2.17%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('E44').Value = '  -21.12%  '
$ws.Range('D45').Value = "'153.02"
$ws.Range('E45').Value = '  +4.33%  '
$ws.Range('E46').Value = '  +1.99%  '
$ws.Range('E47').Value = '  -1.81%  '
$ws.Range('D48').Value = "'20.32"
$ws.Range('E48').Value = '  -1.72%  '
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('D50').Value = "'0.0506"
$ws.Range('E50').Value = '  -1.80%  '
$ws.Range('D51').Value = "'0.0910"
$ws.Range('E51').Value = '  -1.26%  '
